$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the leftover "fill" formatting that was applied to row 3 (id 2 / The HALLOWEEN Chapter)
$ws.Rows(3).ClearFormats()

# Add new DLC chapter row (row 48): Sinister Grace
$ws.Range("A48").Value = 60
$ws.Range("B48").Value = "Sinister Grace"
$ws.Range("C48").Value = "23.09.2025"
$ws.Range("D48").Value = 1
$ws.Range("E48").Formula = "=CHOOSE(D48, ""Chapter DLC"", ""Half-Chapter DLC"", ""Clothing Pack DLC"", ""Original Soundtrack DLC"", ""Character Pack DLC"", ""Other"", ""Retracted"", ""Chapter Pack DLC"")"
$ws.Range("F48").Value = "9.2.0"
$ws.Range("G48").Value = "Maple"
$ws.Range("H48").Formula = "=G48"
$ws.Range("I48").Value = "Vee Boonyasak; Krasue"

# Start of a new row (row 49) - only the id has been typed in so far
$ws.Range("A49").Value = 61

# Reflect where the user had scrolled to / what was selected when the file was saved
$ws.Range("B49").Select()
